$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert the 4 new rows first (top to bottom), using the row numbers they
#    will occupy in the FINAL layout. Each insert pushes everything at/after
#    that row down by one, so later insert targets already account for the
#    earlier inserts.
# ---------------------------------------------------------------------------
$ws.Rows.Item(6).Insert()   # new "Killer(R) ... 22.250.0.4" row in Bad Drivers
$ws.Rows.Item(6).Insert()   # new "Intel(R) ... 23.20.1.1" row in Bad Drivers
$ws.Rows.Item(17).Insert()  # new "Intel(R) ... 21.40.1.3" row in Good Drivers
$ws.Rows.Item(19).Insert()  # new "Killer(R) ... 22.250.0.4" row in Good Drivers

# ---------------------------------------------------------------------------
# 2) Widen column A.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 80.17

# ---------------------------------------------------------------------------
# 3) Bad Drivers table updates (rows 3-9 after the inserts above).
# ---------------------------------------------------------------------------
$ws.Range("C3").Value = 1939
$ws.Range("D3").Value = 65.8

$ws.Range("B4").Value = 13
$ws.Range("C4").Value = 3969
$ws.Range("D4").Value = 92.8

$ws.Range("D5").Value = 95.8

# New row 6
$ws.Range("A6").Value = "Killer(R) Wi-Fi 6 AX1650s 160MHz Wireless Network Adapter (201D2W) - 22.250.0.4"
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 8
$ws.Range("D6").Value = 97.2

# New row 7
$ws.Range("A7").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.20.1.1"
$ws.Range("B7").Value = 4
$ws.Range("C7").Value = 486
$ws.Range("D7").Value = 97.4

# Row 8 (was old row 6 - AX211 23.30.0.6)
$ws.Range("C8").Value = 11

# Row 9 (Totals)
$ws.Range("B9").Value = 35
$ws.Range("C9").Value = 6752

# ---------------------------------------------------------------------------
# 4) Good Drivers table updates (rows 15-27 after the inserts above).
# ---------------------------------------------------------------------------

# New row 17 - give B17 the same "Total Samples" number style as the data
# rows below it (style was inherited from the header row on Insert()).
$ws.Range("B18").Copy()
$ws.Range("B17").PasteSpecial(-4122)
$ws.Range("A17").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.40.1.3"
$ws.Range("B17").Value = 11128
$ws.Range("D17").Value = 100

# Row 18 (was old row 15 - AX201 23.100.0.4)
$ws.Range("B18").Value = 486214

# New row 19 - "Driver Vintage" text needs to stay literal text, not get
# auto-converted to a date serial. Force text format, enter it, then copy
# the correct number/border format back on top (value itself is untouched
# by a formats-only paste).
$ws.Range("A19").Value = "Killer(R) Wi-Fi 6 AX1650s 160MHz Wireless Network Adapter (201D2W) - 22.250.0.4"
$ws.Range("B19").Value = 58842
$ws.Range("D19").Value = 99.9
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2023-07-25"
$ws.Range("D18").Copy()
$ws.Range("E19").PasteSpecial(-4122)

# Row 20 (was old row 16 - AX211 22.150.3.1)
$ws.Range("B20").Value = 11140

# Row 21 (was old row 17 - AX211 22.150.0.3)
$ws.Range("B21").Value = 14487

# Row 22 (was old row 18 - AX211 22.100.1.1) - unchanged

# Row 23 (was old row 19 - AX201 22.80.0.9)
$ws.Range("B23").Value = 79953

# Row 24 (was old row 20 - AX201 22.50.1.1)
$ws.Range("B24").Value = 35355

# Row 25 (was old row 21 - AX201 21.110.3.2)
$ws.Range("B25").Value = 65425

# Row 26 (was old row 22 - AX201 21.70.0.6)
$ws.Range("B26").Value = 117653

# Row 27 (was old row 23 - AX201 21.60.2.1) - unchanged
